$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999138966,
    0.99461792622087897,
    0.97308323888104287,
    0.9630690174249974,
    0.95345794703980391,
    0.93183366137794166,
    0.92804480965330405,
    0.92359053460479523,
    0.91966820336726185,
    0.91667788066554401,
    0.9163490354348599,
    0.91597931064069082,
    0.91591687795053367,
    0.91664345554294036,
    0.91405209972053525,
    0.91154569030009913,
    0.90783786577374725,
    0.90672897042652278,
    0.99288333742582802,
    0.98576641247042729,
    0.98436794913589165,
    0.98310344817936968,
    0.96604767738514175,
    0.95302601239174578,
    0.94656889285287016,
    0.92915917352868349,
    0.92567914406889673,
    0.91299935598866566,
    0.90372103206933563,
    0.89927351066613603,
    0.89162012403765312,
    0.88994084627163517,
    0.88942084791301745
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
